$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metrics por Modelo" ---
$ws1 = $wb.Worksheets.Item("Metrics por Modelo")

$ws1.Range("B2").Value = 0.8152866242038217
$ws1.Range("C2").Value = 0.768
$ws1.Range("D2").Value = 0.75
$ws1.Range("E2").Value = 0.7868852459016393

$ws1.Range("B3").Value = 0.8089171974522293
$ws1.Range("C3").Value = 0.7368421052631579
$ws1.Range("D3").Value = 0.65625
$ws1.Range("E3").Value = 0.84

$ws1.Range("B4").Value = 0.8089171974522293
$ws1.Range("C4").Value = 0.7413793103448276
$ws1.Range("D4").Value = 0.671875
$ws1.Range("E4").Value = 0.8269230769230769

# --- Sheet 2: "Hiperparam + Rendimiento" ---
$ws2 = $wb.Worksheets.Item("Hiperparam + Rendimiento")

$ws2.Range("F2").Value = 0.7906976744186046
$ws2.Range("G2").Value = 0.8343949044585988
$ws2.Range("H2").Value = 0.8152866242038217

$ws2.Range("F3").Value = 0.8344733242134063
$ws2.Range("G3").Value = 0.8152866242038217
$ws2.Range("H3").Value = 0.8089171974522293

$ws2.Range("F4").Value = 0.8331053351573188
$ws2.Range("G4").Value = 0.8343949044585988
$ws2.Range("H4").Value = 0.8089171974522293
